# Complete Hawaii (HI) scraper - 13 staff members
# Also tidies up the GA sheet (clears stray empty trailing cells picked up by
# the earlier scrape pass, autosizes its two content columns) and leaves the
# HI tab as the active / selected sheet in the workbook.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# GA: drop the empty C:F placeholder cells that trailed every staff
# row, then resize columns A/B to fit their (now sole) content.
# -----------------------------------------------------------------
$ga = $wb.Worksheets.Item("GA")
$ga.Range("C4:F28").ClearContents()
$ga.Columns.Item(1).ColumnWidth = 16.09
$ga.Columns.Item(2).ColumnWidth = 50.59
$ga.Range("C30").Select()

# -----------------------------------------------------------------
# HI: append the 13 INNOVATE Hawaii staff members scraped from
# https://www.htdc.org/our-team/
# -----------------------------------------------------------------
$hi = $wb.Worksheets.Item("HI")

$hi.Range("A4").Value = "Trung Lam"
$hi.Range("B4").Value = "Executive Director & CEO / Acting MEP Center Director"
$hi.Range("C4:F4").Borders.LineStyle = 0

$hi.Range("A5").Value = "Matthew Kobayashi"
$hi.Range("B5").Value = "Project Development Manager"
$hi.Range("C5:F5").Borders.LineStyle = 0

$hi.Range("A6").Value = "Sandi Kanemori"
$hi.Range("B6").Value = "Sr. Economic Program Manager"
$hi.Range("C6:F6").Borders.LineStyle = 0

$hi.Range("A7").Value = "Umma Berkelman"
$hi.Range("B7").Value = "Economic Development Specialist"
$hi.Range("C7:F7").Borders.LineStyle = 0

$hi.Range("A8").Value = "Cindy Matsuki"
$hi.Range("B8").Value = "Economic Development Specialist - HSBIR"
$hi.Range("C8:F8").Borders.LineStyle = 0

$hi.Range("A9").Value = "Karlton Tomomitsu"
$hi.Range("B9").Value = "Economic Development Specialist"
$hi.Range("C9:F9").Borders.LineStyle = 0

$hi.Range("A10").Value = "Dave Molinaro"
$hi.Range("B10").Value = "HCATT Director"
$hi.Range("C10:F10").Borders.LineStyle = 0

$hi.Range("A11").Value = "Kristy Carpio"
$hi.Range("B11").Value = "HCATT Project Manager"
$hi.Range("C11:F11").Borders.LineStyle = 0

$hi.Range("A12").Value = "Tuan La"
$hi.Range("B12").Value = "HI-CAP Manager"
$hi.Range("C12:F12").Borders.LineStyle = 0

$hi.Range("A13").Value = "Wayne Layugan"
$hi.Range("B13").Value = "Program Manager"
$hi.Range("C13:F13").Borders.LineStyle = 0

$hi.Range("A14").Value = "Wendy Oshiro"
$hi.Range("B14").Value = "Project Manager"
$hi.Range("C14:F14").Borders.LineStyle = 0

$hi.Range("A15").Value = "Ray Gomez"
$hi.Range("B15").Value = "Chief Financial Officer"
$hi.Range("C15:F15").Borders.LineStyle = 0

$hi.Range("A16").Value = "Stephanie Yuu-Sato"
$hi.Range("B16").Value = "Contracts & Project Manager"
$hi.Range("C16:F16").Borders.LineStyle = 0

# HI becomes the active sheet / selected tab, cursor back at A1.
$hi.Range("A1").Select()
$hi.Activate()

Write-Host "HI staff rows added; GA trimmed."
